# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.983.43'
$ws.Range("E2").Value = '''  +0.78%  '
$ws.Range("D3").Value = '''2.570.09'
$ws.Range("E3").Value = '''  -0.32%  '
$ws.Range("E4").Value = '''  +0.07%  '
$ws.Range("D5").Value = '''565.21'
$ws.Range("E5").Value = '''  +4.40%  '
$ws.Range("D6").Value = '''142.29'
$ws.Range("E6").Value = '''  -1.32%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '''  +0.05%  '
$ws.Range("D8").Value = '''0.592'
$ws.Range("E8").Value = '''  +1.50%  '
$ws.Range("D9").Value = '''2.576.25'
$ws.Range("E9").Value = '''  -0.17%  '
$ws.Range("D10").Value = '''6.64'
$ws.Range("E10").Value = '''  -1.88%  '
$ws.Range("D11").Value = '''0.102'
$ws.Range("E11").Value = '''  +2.17%  '
$ws.Range("D12").Value = '''0.152'
$ws.Range("E12").Value = '''  +9.50%  '
$ws.Range("D13").Value = '''0.340'
$ws.Range("E13").Value = '''  +2.04%  '
$ws.Range("D14").Value = '''3.026.65'
$ws.Range("E14").Value = '''  -0.18%  '
$ws.Range("D15").Value = '''59.095.23'
$ws.Range("E15").Value = '''  +1.16%  '
$ws.Range("D16").Value = '''21.74'
$ws.Range("E16").Value = '''  +5.79%  '
$ws.Range("D17").Value = '''0.0000135'
$ws.Range("E17").Value = '''  +3.03%  '
$ws.Range("D18").Value = '''2.576.40'
$ws.Range("E18").Value = '''  +1.20%  '
$ws.Range("D19").Value = '''4.49'
$ws.Range("E19").Value = '''  +0.54%  '
$ws.Range("D20").Value = '''334.40'
$ws.Range("E20").Value = '''  +0.11%  '
$ws.Range("D21").Value = '''10.13'
$ws.Range("E21").Value = '''  +0.88%  '
$ws.Range("D22").Value = '''6.16'
$ws.Range("E22").Value = '''  +1.17%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '''  +0.14%  '
$ws.Range("D24").Value = '''64.69'
$ws.Range("E24").Value = '''  -2.50%  '
$ws.Range("D25").Value = '''0.445'
$ws.Range("E25").Value = '''  +5.45%  '
$ws.Range("E26").Value = '''  +0.47%  '
$ws.Range("E27").Value = '''  +1.63%  '
$ws.Range("D28").Value = '''7.20'
$ws.Range("E28").Value = '''  +2.12%  '
$ws.Range("D29").Value = '''0.0₃0778'
$ws.Range("E29").Value = '''  +5.18%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '''  -0.01%  '
$ws.Range("D31").Value = '''1.68'
$ws.Range("E31").Value = '''  +2.53%  '
$ws.Range("D32").Value = '''160.20'
$ws.Range("E32").Value = '''  +4.76%  '
$ws.Range("D33").Value = '''6.02'
$ws.Range("E33").Value = '''  +0.32%  '
$ws.Range("D34").Value = '''18.83'
$ws.Range("E34").Value = '''  -0.48%  '
$ws.Range("D35").Value = '''4.00'
$ws.Range("E35").Value = '''  +2.49%  '
$ws.Range("D36").Value = '''0.877'
$ws.Range("E36").Value = '''  +3.09%  '
$ws.Range("D37").Value = '''0.878'
$ws.Range("E37").Value = '''  +6.98%  '
$ws.Range("D38").Value = '''1.12'
$ws.Range("E38").Value = '''  +2.85%  '
$ws.Range("E39").Value = '''  -0.85%  '
$ws.Range("D40").Value = '''1.48'
$ws.Range("E40").Value = '''  +4.36%  '
$ws.Range("D41").Value = '''294.68'
$ws.Range("E41").Value = '''  +5.85%  '
$ws.Range("D42").Value = '''3.62'
$ws.Range("E42").Value = '''  +0.98%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("D44").Value = '''0.0971'
$ws.Range("E44").Value = '''  +2.89%  '
$ws.Range("B45").Value = '''Mantle'
$ws.Range("C45").Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.590'
$ws.Range("E45").Value = '''  -0.08%  '
$ws.Range("B46").Value = '''Hedera'
$ws.Range("C46").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0534'
$ws.Range("E46").Value = '''  +0.89%  '
$ws.Range("B47").Value = '''WhiteBITCoin'
$ws.Range("C47").Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '''10.62'
$ws.Range("E47").Value = '''  -0.03%  '
$ws.Range("D48").Value = '''124.60'
$ws.Range("E48").Value = '''  +14.07%  '
$ws.Range("D49").Value = '''18.91'
$ws.Range("E49").Value = '''  +2.24%  '
$ws.Range("D50").Value = '''0.0230'
$ws.Range("E50").Value = '''  +2.01%  '
$ws.Range("B51").Value = '''InjectiveProtocol'
$ws.Range("C51").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '''18.31'
$ws.Range("E51").Value = '''  +2.54%  '
